# BIS-769: Fixed xls test files
#
# Adds two new trailing columns ("Pattern" / "Pattern Type") to the two
# property-definition header rows (row 4 and row 12) of the sample-type
# export sheet, mirroring the style of the existing "Unique" header cell
# (column L) onto the two new columns (M, N).
#
# xlPasteFormats = -4122 : pastes only the cell format (number format,
# font, fill, borders, alignment, protection) without touching the
# clipboard source's value, so the destination cell keeps reusing the
# very same shared style record (s="6") instead of minting a new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 4 header ---------------------------------------------------
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial($xlPasteFormats)
$ws.Range("M4").Value = "Pattern"

$ws.Range("L4").Copy()
$ws.Range("N4").PasteSpecial($xlPasteFormats)
$ws.Range("N4").Value = "Pattern Type"

# --- Row 12 header (second SAMPLE_TYPE block) ------------------------
$ws.Range("L12").Copy()
$ws.Range("M12").PasteSpecial($xlPasteFormats)
$ws.Range("M12").Value = "Pattern"

$ws.Range("L12").Copy()
$ws.Range("N12").PasteSpecial($xlPasteFormats)
$ws.Range("N12").Value = "Pattern Type"

# Clear the marching-ants clipboard marquee and match the author's
# final selection (M12:N12, active cell M12), as recorded in the sheet.
$excel.CutCopyMode = $false
[void]$ws.Range("M12:N12").Select()
